$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Bump the quoted date by one day (2024-01-19 -> 2024-01-20)
$ws.Range("A1").Value = 45311

# Reduce the listed prices for the three "Bisagra t/ ESCALERA" items
$ws.Range("D26").Value = 175.344
$ws.Range("D27").Value = 221.602
$ws.Range("D28").Value = 300
